# Rename the output-commodity trade-link codes from "TB_H2_..." to
# "TB_H2GC_..." on the SUP_TRADE sheet.
#
# Each of these 8 names lives once in the shared-string table, anchored at
# the "source" cell in column O for each 4-row trade-link block (O4, O9,
# O14, O19, O24, O29, O34, O39). Every other O cell in a block is a live
# formula ("=O<prev row>") that simply copies the value down, so updating
# the anchor cell's text automatically ripples the new text (and the
# cached <v> shown by each formula) through the rest of the block once
# Excel recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUP_TRADE")

$ws.Range("O4").Value  = "TB_H2GC_DKISLBH_DKE_01"
$ws.Range("O9").Value  = "TB_H2GC_DKISLBH_DKE_02"
$ws.Range("O14").Value = "TB_H2GC_DKISL1_DKW_01"
$ws.Range("O19").Value = "TB_H2GC_DKISL1_DKW_02"
$ws.Range("O24").Value = "TB_H2GC_DKISL2_DKW_01"
$ws.Range("O29").Value = "TB_H2GC_DKISL2_DKW_02"
$ws.Range("O34").Value = "TB_H2GC_DKISL3_DKW_01"
$ws.Range("O39").Value = "TB_H2GC_DKISL3_DKW_02"
